# Applies the "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
# edit to the NIT-9012768467 estado de cuenta workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header totals: Valor Mora total and Cant. Periodos increase (new worker
#    Brenda adds one more period, and Lida now has 6 periods instead of 5).
# ---------------------------------------------------------------------------
$ws.Range("E11").Value2 = 393640
$ws.Range("F13").Value2 = 7

# ---------------------------------------------------------------------------
# 2) Prepare row 22 (new last data row) by copying the formatting that
#    currently "closes" the table on row 21 (thicker bottom border etc.)
# ---------------------------------------------------------------------------
$ws.Range("B21:J21").Copy()
$ws.Range("B22:J22").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# Now make row 21 use the regular "middle of table" formatting (same as rows
# 16-20) since it will no longer be the last row.
$ws.Range("B16:J16").Copy()
$ws.Range("B21:J21").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Rewrite the worker/period data block (rows 16-22).
#    Row 16: new worker Brenda Carolina Marimon Marimon (period 2411)
#    Rows 17-22: Lida Rosa Amaranto Correa, periods 2503-2508 (ascending,
#    with the brand-new period 2508 appended).
# ---------------------------------------------------------------------------
$ws.Range("B16").Value2 = "CC"
$ws.Range("C16").Value2 = "1047485132"
$ws.Range("D16").Value2 = "BRENDA CAROLINA MARIMON MARIMON"
$ws.Range("E16").Value2 = "2411"
$ws.Range("F16").Value2 = 52000
$ws.Range("G16").Value2 = 1423500

$periods = @("2503", "2504", "2505", "2506", "2507", "2508")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 17 + $i
    $ws.Range("B$r").Value2 = "CC"
    $ws.Range("C$r").Value2 = "45462518"
    $ws.Range("D$r").Value2 = "LIDA ROSA AMARANTO CORREA"
    $ws.Range("E$r").Value2 = $periods[$i]
    $ws.Range("F$r").Value2 = 56940
    $ws.Range("G$r").Value2 = 1423500
}

# ---------------------------------------------------------------------------
# 4) Move the signature block down one row (row 26 -> 27, row 27 -> 28) to
#    make room for the extra data row above.
# ---------------------------------------------------------------------------
$ws.Range("B26:J26").Copy()
$ws.Range("B27:J27").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("B27:J27").Copy()
$ws.Range("B28:J28").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("B27").Value2 = "___________________________________"
$ws.Range("H27").Value2 = "___________________________________"
$ws.Range("B26").ClearContents()
$ws.Range("H26").ClearContents()

$ws.Range("B28").Value2 = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H28").Value2 = "FIRMA DEL REPRESENTANTE LEGAL"

$ws.Range("B28:C28").Merge()
$ws.Range("H28:J28").Merge()
$ws.Range("B26:C26").UnMerge()
$ws.Range("H26:J26").UnMerge()

Write-Host "Edit applied"
